# Update the Cupid slot review copy: new title/meta wording, and refreshed
# "What we like" / "What we don't like" bullet points.
#
# Each change is a like-for-like text swap, so a scoped Find/Replace across
# the whole document body is the natural Word automation approach - it
# preserves paragraph styles (Heading1/2, ListBullet) and run formatting
# (bold/italic) untouched, only rewriting the literal wording.

$d = $word.ActiveDocument

function Replace-Text {
    param([string]$old, [string]$new)

    $found = $d.Content.Find.Execute(
        $old,    # FindText
        $true,   # MatchCase
        $false,  # MatchWholeWord
        $false,  # MatchWildcards
        $false,  # MatchSoundsLike
        $false,  # MatchAllWordForms
        $true,   # Forward
        1,       # Wrap (wdFindContinue)
        $false,  # Format
        $new,    # ReplaceWith
        2        # Replace (wdReplaceAll)
    )
    Write-Output "Replace '$old' -> '$new': $found"
}

# Page title / heading (also reused verbatim later as the bold "title" line).
Replace-Text "Play Cupid Slot Game for Free | Review 2021" "Play Cupid Slot for Free - Review and Gameplay"

# "What we like" bullets.
Replace-Text "Beautiful graphics and immersive arcade music" "Well-defined graphics and appealing theme"
Replace-Text "Double winnings with the Wild symbol" "Large grid and easy-to-use interface"
Replace-Text "Free spins with potential for high multipliers" "Option to autoplay with customizable limits"
Replace-Text "Gamble feature allows for doubling of winnings" "Gamble function for extra excitement"

# "What we don't like" bullets.
Replace-Text "No progressive jackpot" "Limited number of paylines"
Replace-Text "Not much variety in terms of gameplay features" "Gamble function can result in loss of winnings"

# Closing meta description (italic run).
Replace-Text "Read our unbiased review of Cupid slot game, play for free and enjoy the beautiful graphics, free spins, and high multipliers of this online slot game." "Read our review of Cupid slot game and play for free. Enjoy the charming graphics and exciting features."

Write-Output "Done"
